$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '29.127.90'
Set-TextValue $ws.Range("E2") '  +0.02%  '

Set-TextValue $ws.Range("D3") '1.900.96'
Set-TextValue $ws.Range("E3") '  -0.22%  '

Set-TextValue $ws.Range("D4") '0.9997'
Set-TextValue $ws.Range("E4") '  -0.29%  '

Set-TextValue $ws.Range("D5") '325.42'
Set-TextValue $ws.Range("E5") '  -0.46%  '

Set-TextValue $ws.Range("D6") '1.000'
Set-TextValue $ws.Range("E6") '  -0.23%  '

Set-TextValue $ws.Range("D7") '0.4625'
Set-TextValue $ws.Range("E7") '  -0.11%  '

Set-TextValue $ws.Range("D8") '0.3902'
Set-TextValue $ws.Range("E8") '  -0.78%  '

Set-TextValue $ws.Range("D9") '0.07879'
Set-TextValue $ws.Range("E9") '  -0.90%  '

Set-TextValue $ws.Range("D10") '0.9916'
Set-TextValue $ws.Range("E10") '  -0.70%  '

Set-TextValue $ws.Range("D11") '21.90'
Set-TextValue $ws.Range("E11") '  -1.42%  '

Set-TextValue $ws.Range("D12") '1.916.35'
Set-TextValue $ws.Range("E12") '  +2.26%  '

Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '5.777'
Set-TextValue $ws.Range("E13") '  +0.25%  '

Set-TextValue $ws.Range("B14") 'Chainlink'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D14") '7.067'
Set-TextValue $ws.Range("E14") '  -0.28%  '

Set-TextValue $ws.Range("D15") '0.06986'
Set-TextValue $ws.Range("E15") '  +0.85%  '

Set-TextValue $ws.Range("D16") '87.96'
Set-TextValue $ws.Range("E16") '  -0.53%  '

Set-TextValue $ws.Range("E17") '  -0.21%  '

Set-TextValue $ws.Range("D18") '0.000009916'
Set-TextValue $ws.Range("E18") '  -1.47%  '

Set-TextValue $ws.Range("E19") '  -0.39%  '

Set-TextValue $ws.Range("E20") '  -0.19%  '

Set-TextValue $ws.Range("D21") '29.143.52'
Set-TextValue $ws.Range("E21") '  +0.01%  '

Set-TextValue $ws.Range("D22") '5.319'
Set-TextValue $ws.Range("E22") '  -0.82%  '

Set-TextValue $ws.Range("D23") '11.11'
Set-TextValue $ws.Range("E23") '  +0.14%  '

Set-TextValue $ws.Range("D24") '2.123.84'
Set-TextValue $ws.Range("E24") '  -0.80%  '

Set-TextValue $ws.Range("D25") '2.110'
Set-TextValue $ws.Range("E25") '  +2.45%  '

Set-TextValue $ws.Range("D26") '155.87'
Set-TextValue $ws.Range("E26") '  -0.45%  '

Set-TextValue $ws.Range("D27") '19.40'
Set-TextValue $ws.Range("E27") '  -0.26%  '

Set-TextValue $ws.Range("D28") '5.919'
Set-TextValue $ws.Range("E28") '  -1.29%  '

Set-TextValue $ws.Range("D29") '118.72'
Set-TextValue $ws.Range("E29") '  -0.27%  '

Set-TextValue $ws.Range("D30") '1.880'
Set-TextValue $ws.Range("E30") '  -5.60%  '

Set-TextValue $ws.Range("D31") '0.09329'
Set-TextValue $ws.Range("E31") '  -0.66%  '

Set-TextValue $ws.Range("D32") '0.9000'
Set-TextValue $ws.Range("E32") '  -2.31%  '

Set-TextValue $ws.Range("D33") '5.251'
Set-TextValue $ws.Range("E33") '  -1.52%  '

Set-TextValue $ws.Range("D34") '1.326'
Set-TextValue $ws.Range("E34") '  -1.67%  '

Set-TextValue $ws.Range("D35") '3.156'
Set-TextValue $ws.Range("E35") '  -3.22%  '

Set-TextValue $ws.Range("D36") '0.05805'
Set-TextValue $ws.Range("E36") '  -0.08%  '

Set-TextValue $ws.Range("E37") '  -2.10%  '

Set-TextValue $ws.Range("D38") '0.02087'
Set-TextValue $ws.Range("E38") '  -0.57%  '

Set-TextValue $ws.Range("D39") '1.000'
Set-TextValue $ws.Range("E39") '  -0.12%  '

Set-TextValue $ws.Range("D40") '7.730'
Set-TextValue $ws.Range("E40") '  -2.71%  '

Set-TextValue $ws.Range("E41") '  -0.80%  '

Set-TextValue $ws.Range("D42") '0.1793'
Set-TextValue $ws.Range("E42") '  -0.12%  '

Set-TextValue $ws.Range("D43") '9.746'
Set-TextValue $ws.Range("E43") '  -2.09%  '

Set-TextValue $ws.Range("B44") 'RenderToken'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D44") '2.243'
Set-TextValue $ws.Range("E44") '  +1.78%  '

Set-TextValue $ws.Range("B45") 'EnergySwap'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '11.89'
Set-TextValue $ws.Range("E45") '  -0.94%  '

Set-TextValue $ws.Range("D46") '0.5361'
Set-TextValue $ws.Range("E46") '  -1.07%  '

Set-TextValue $ws.Range("D47") '0.07008'

Set-TextValue $ws.Range("D49") '2.552'
Set-TextValue $ws.Range("E49") '  +0.27%  '

Set-TextValue $ws.Range("D50") '113.10'
Set-TextValue $ws.Range("E50") '  +0.96%  '

Set-TextValue $ws.Range("D51") '1.042'
Set-TextValue $ws.Range("E51") '  -0.77%  '
